# The workbook gained one new data row (a new daily price record), inserted
# at row 12, pushing every subsequent row down by one. This reproduces that
# insertion and fills the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; this shifts rows 12..115 down to 13..116
# and carries the existing formatting (e.g. the date style on column D) down
# with it / into the newly inserted row.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new record's data. All of the
# non-price/date columns are identical to the (now shifted) row below it.
$ws.Cells.Item(12, 1).Value  = 10
$ws.Cells.Item(12, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(12, 3).Value  = "La Araucanía"
$ws.Cells.Item(12, 4).Value  = 45111
$ws.Cells.Item(12, 5).Value  = 9
$ws.Cells.Item(12, 6).Value  = 300000001
$ws.Cells.Item(12, 7).Value  = "Rabanito"
$ws.Cells.Item(12, 8).Value  = "Sin especificar"
$ws.Cells.Item(12, 9).Value  = "Primera"
$ws.Cells.Item(12, 10).Value = 30
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 8000
$ws.Cells.Item(12, 13).Value = 8000
$ws.Cells.Item(12, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(12, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(12, 16).Value = 667
$ws.Cells.Item(12, 17).Value = 12
$ws.Cells.Item(12, 18).Value = "Hortaliza"
